$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/29/2024  Through  2/4/2024"

# --- Main crime-stat table updates (rows 14-30) ---
$ws.Range("C14").Value = "'0"
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = -50
$ws.Range("N14").Value = -86.666666666666
$ws.Range("C15").Value = 3
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -42.105263157894
$ws.Range("L15").Value = -42.105263157894
$ws.Range("M15").Value = -21.428571428571
$ws.Range("N15").Value = -47.619047619047
$ws.Range("C16").Value = 45
$ws.Range("D16").Value = 47
$ws.Range("E16").Value = -4.255319148936
$ws.Range("F16").Value = 161
$ws.Range("G16").Value = 143
$ws.Range("H16").Value = 12.587412587412
$ws.Range("I16").Value = 197
$ws.Range("J16").Value = 176
$ws.Range("K16").Value = 11.931818181818
$ws.Range("L16").Value = 60.162601626016
$ws.Range("M16").Value = -12.053571428571
$ws.Range("N16").Value = -81.519699812382
$ws.Range("C17").Value = 56
$ws.Range("D17").Value = 44
$ws.Range("E17").Value = 27.272727272727
$ws.Range("F17").Value = 194
$ws.Range("G17").Value = 174
$ws.Range("H17").Value = 11.494252873563
$ws.Range("I17").Value = 255
$ws.Range("J17").Value = 233
$ws.Range("K17").Value = 9.442060085836
$ws.Range("L17").Value = 21.428571428571
$ws.Range("M17").Value = 94.656488549618
$ws.Range("N17").Value = 3.658536585365
$ws.Range("C18").Value = 38
$ws.Range("E18").Value = -15.555555555555
$ws.Range("F18").Value = 149
$ws.Range("G18").Value = 181
$ws.Range("H18").Value = -17.679558011049
$ws.Range("I18").Value = 185
$ws.Range("J18").Value = 217
$ws.Range("K18").Value = -14.746543778801
$ws.Range("L18").Value = 1.092896174863
$ws.Range("M18").Value = -35.540069686411
$ws.Range("N18").Value = -88.915518274415
$ws.Range("C19").Value = 145
$ws.Range("D19").Value = 144
$ws.Range("E19").Value = 0.694444444444
$ws.Range("F19").Value = 539
$ws.Range("G19").Value = 501
$ws.Range("H19").Value = 7.584830339321
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = 616
$ws.Range("K19").Value = 5.519480519480
$ws.Range("L19").Value = -22.341696535244
$ws.Range("M19").Value = 63.727959697733
$ws.Range("N19").Value = -10.958904109589
$ws.Range("C20").Value = 37
$ws.Range("D20").Value = 45
$ws.Range("E20").Value = -17.777777777777
$ws.Range("F20").Value = 166
$ws.Range("G20").Value = 169
$ws.Range("H20").Value = -1.775147928994
$ws.Range("I20").Value = 212
$ws.Range("J20").Value = 204
$ws.Range("K20").Value = 3.921568627450
$ws.Range("L20").Value = 53.623188405797
$ws.Range("M20").Value = 30.864197530864
$ws.Range("N20").Value = -91.630477694433
$ws.Range("C21").Value = 324
$ws.Range("D21").Value = 328
$ws.Range("E21").Value = -1.219512195121
$ws.Range("F21").Value = 1219
$ws.Range("G21").Value = 1182
$ws.Range("H21").Value = 3.130287648054
$ws.Range("I21").Value = 1512
$ws.Range("J21").Value = 1467
$ws.Range("K21").Value = 3.067484662576
$ws.Range("L21").Value = -0.066093853271
$ws.Range("M21").Value = 24.239934264585
$ws.Range("N21").Value = -75.923566878980
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -20
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = 27.777777777777
$ws.Range("I22").Value = 25
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 19.047619047619
$ws.Range("L22").Value = -7.407407407407
$ws.Range("M22").Value = 4.166666666666
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 6
$ws.Range("H23").Value = -68.421052631578
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = -27.272727272727
$ws.Range("L23").Value = -30.434782608695
$ws.Range("M23").Value = 14.285714285714
$ws.Range("C24").Value = 361
$ws.Range("D24").Value = 280
$ws.Range("E24").Value = 28.928571428571
$ws.Range("F24").Value = 1331
$ws.Range("G24").Value = 1163
$ws.Range("H24").Value = 14.445399828031
$ws.Range("I24").Value = 1635
$ws.Range("J24").Value = 1446
$ws.Range("K24").Value = 13.070539419087
$ws.Range("L24").Value = 22.747747747747
$ws.Range("M24").Value = 95.107398568019
$ws.Range("D25").Value = 84
$ws.Range("E25").Value = 13.095238095238
$ws.Range("F25").Value = 364
$ws.Range("G25").Value = 358
$ws.Range("H25").Value = 1.675977653631
$ws.Range("I25").Value = 459
$ws.Range("J25").Value = 462
$ws.Range("K25").Value = -0.649350649350
$ws.Range("L25").Value = 11.678832116788
$ws.Range("M25").Value = 10.336538461538
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = -5.555555555555
$ws.Range("I26").Value = 21
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -22.222222222222
$ws.Range("L26").Value = -19.230769230769
$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 30
$ws.Range("F27").Value = 42
$ws.Range("G27").Value = 48
$ws.Range("H27").Value = -12.5
$ws.Range("I27").Value = 49
$ws.Range("J27").Value = 55
$ws.Range("K27").Value = -10.909090909090
$ws.Range("L27").Value = 16.666666666666
$ws.Range("C28").Value = "'0"
$ws.Range("N28").Value = -93.333333333333
$ws.Range("C29").Value = "'0"
$ws.Range("N29").Value = -93.333333333333
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 6
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = 50
$ws.Range("L30").Value = 100
